$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the table-name typos so both rows refer to the same table name
$ws.Range("B16").Value = "Tabela vendas_itens"
$ws.Range("B17").Value = "Tabela vendas_itens"

# Fix the typo / wording in the note cell (A2)
$ws.Range("A2").Value = "Obs.: Todas as tabelas e funções serão em aquivos individuais"

# Clear the "!" status markers that were left over on these Insert rows
$ws.Range("F7").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("F15").Value = ""

# Update the selection shown when the sheet is reopened
$ws.Range("A1:E1").Select()
